$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Cells.Item(6, 2).Value = "active"

# Date: updated publish timestamp
$ws.Cells.Item(8, 2).Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: (was reusing "false") -> "true"
# Use a leading apostrophe so the literal "true" is stored as text rather
# than being auto-converted to a boolean, then restore the original cell
# formatting (which the apostrophe entry perturbs) by copying the format
# from a neighboring cell that already carries the unchanged style.
$ws.Cells.Item(17, 2).Value = "'true"
$ws.Cells.Item(18, 2).Copy()
$ws.Cells.Item(17, 2).PasteSpecial(-4122)
$excel.CutCopyMode = 0
